# The roboticRNAPrep column (H) was stored as Boolean cells (custom
# "TRUE"/"FALSE" number format). The author retyped these as the literal
# text word "False" instead of the boolean FALSE, so the column is
# switched to a Text-formatted column of string cells reading "False".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A plain `Range.Value = "False"` assignment is auto-coerced back into a
# Boolean by Excel's type inference (same as typing FALSE straight into a
# cell), so build the literal string via a text formula in a scratch cell
# and paste its *value* across - a value-only paste is not re-inferred.
$scratch = $ws.Cells.Item(1, 26)
for ($r = 2; $r -le 27; $r++) {
    $scratch.Formula = '="False"'
    $scratch.Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4163)  # xlPasteValues
}
$scratch.ClearContents()
$excel.CutCopyMode = $false

# Re-format the column as text (was the custom boolean TRUE/FALSE format).
$ws.Range("H2:H27").NumberFormat = "@"

# Match the selection left behind in the saved file.
$ws.Range("H2:H27").Select()
